$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row for the "Invalid Click 2" interaction.
$ws.Range("A13").Value = "Invalid Click 2"
$ws.Range("B13").Value = "invalid_click_2.mp3"
$ws.Range("C13").Value = "figma"

# Update the "背景人声" row's file name to the re-versioned filenames.
$ws.Range("B7").Value = "voice_1.mp3//voice_2.mp3"

# Update selection to match the new active cell.
$ws.Range("E14").Select()
